$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# atualizei dados bibi e add
$ws.Range("G3").Value = -52
$ws.Range("G4").Value = 28
$ws.Range("I4").Value = 0.24
$ws.Range("G6").Value = -3
